$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.46
$ws.Range("C2").Value = 0.72
$ws.Range("D2").Value = 26.98
$ws.Range("F2").Value = 7.53
